$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 143.67
$ws.Range("I15").Value = 143.67
$ws.Range("K15").Value = 431.01
$ws.Range("M15").Value = -262.01
$ws.Range("H41").Value = 1110.8462
$ws.Range("I41").Value = 1995
$ws.Range("J41").Value = 717.8889
$ws.Range("K41").Value = 1995
$ws.Range("L41").Value = 717.8889
$ws.Range("M41").Value = -1555
$ws.Range("N41").Value = -1597.8889
$ws.Range("H62").Value = 3017.1765
$ws.Range("I62").Value = 2945.077
$ws.Range("J62").Value = 3251.5
$ws.Range("K62").Value = 2945.077
$ws.Range("L62").Value = 3251.5
$ws.Range("M62").Value = -2321.077
$ws.Range("N62").Value = -4499.5
$ws.Range("H65").Value = 3017.1765
$ws.Range("I65").Value = 2945.077
$ws.Range("J65").Value = 3251.5
$ws.Range("K65").Value = 14725.385
$ws.Range("L65").Value = 16257.5
$ws.Range("M65").Value = -11605.385
$ws.Range("N65").Value = -22497.5
$ws.Range("H86").Value = 3346.8
$ws.Range("I86").Value = 2808.5
$ws.Range("K86").Value = 2808.5
$ws.Range("M86").Value = -1685.5
$ws.Range("H89").Value = 3346.8
$ws.Range("I89").Value = 2808.5
$ws.Range("K89").Value = 14042.5
$ws.Range("M89").Value = -8426.5
$ws.Range("H116").Value = 2221.6667
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 1999
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 1999
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -8883
$ws.Range("H132").Value = 4315322
$ws.Range("I132").Value = 4812745.5
$ws.Range("K132").Value = 14438236.5
$ws.Range("M132").Value = -14435706.5
$ws.Range("H137").Value = 1431.258
$ws.Range("I137").Value = 1330.85
$ws.Range("K137").Value = 3992.55
$ws.Range("M137").Value = -1442.55
$ws.Range("H141").Value = 4446.2
$ws.Range("I141").Value = 4385
$ws.Range("J141").Value = 4589
$ws.Range("K141").Value = 13155
$ws.Range("L141").Value = 13767
$ws.Range("M141").Value = -7975
$ws.Range("N141").Value = -24127

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31970.527
$ws.Range("I32").Value = 14353.8
$ws.Range("J32").Value = 160092.19
$ws.Range("K32").Value = 14353.8
$ws.Range("L32").Value = 160092.19
$ws.Range("M32").Value = -14066.8
$ws.Range("N32").Value = -160666.19
$ws.Range("H44").Value = 12846.125
$ws.Range("J44").Value = 12824.143
$ws.Range("L44").Value = 12824.143
$ws.Range("N44").Value = -13800.143
$ws.Range("H55").Value = 12128.571
$ws.Range("J55").Value = 12128.571
$ws.Range("L55").Value = 12128.571
$ws.Range("N55").Value = -12758.571
$ws.Range("H132").Value = 9965.725
$ws.Range("I132").Value = 10912.099
$ws.Range("J132").Value = 2749.625
$ws.Range("K132").Value = 32736.297
$ws.Range("L132").Value = 8248.875
$ws.Range("M132").Value = -30206.297
$ws.Range("N132").Value = -13308.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H20").Value = 44427.207
$ws.Range("I20").Value = 55427
$ws.Range("K20").Value = 55427
$ws.Range("M20").Value = -55180
$ws.Range("H22").Value = 347.55554
$ws.Range("I22").Value = 335
$ws.Range("K22").Value = 335
$ws.Range("M22").Value = -162
$ws.Range("H60").Value = 44930
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 44930
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 44930
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -46128
$ws.Range("H86").Value = 81480.86
$ws.Range("I86").Value = 113550.6
$ws.Range("J86").Value = 1306.5
$ws.Range("K86").Value = 113550.6
$ws.Range("L86").Value = 1306.5
$ws.Range("M86").Value = -112427.6
$ws.Range("N86").Value = -3552.5
$ws.Range("H89").Value = 81480.86
$ws.Range("I89").Value = 113550.6
$ws.Range("J89").Value = 1306.5
$ws.Range("K89").Value = 567753
$ws.Range("L89").Value = 6532.5
$ws.Range("M89").Value = -562137
$ws.Range("N89").Value = -17764.5
$ws.Range("H99").Value = 1739.9524
$ws.Range("I99").Value = 1499.9286
$ws.Range("J99").Value = 2220
$ws.Range("K99").Value = 1499.9286
$ws.Range("L99").Value = 2220
$ws.Range("M99").Value = -1.92859999999996
$ws.Range("N99").Value = -5216
$ws.Range("H105").Value = 112995
$ws.Range("I105").Value = 78936.84
$ws.Range("J105").Value = 201546.2
$ws.Range("K105").Value = 78936.84
$ws.Range("L105").Value = 201546.2
$ws.Range("M105").Value = -77189.84
$ws.Range("N105").Value = -205040.2
$ws.Range("H134").Value = 3626.6892
$ws.Range("I134").Value = 3766.5667
$ws.Range("J134").Value = 3027.2144
$ws.Range("K134").Value = 11299.7001
$ws.Range("L134").Value = 9081.643199999999
$ws.Range("M134").Value = -8764.7001
$ws.Range("N134").Value = -14151.6432

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 935.8421
$ws.Range("I16").Value = 714
$ws.Range("J16").Value = 1316.1428
$ws.Range("K16").Value = 714
$ws.Range("L16").Value = 1316.1428
$ws.Range("M16").Value = -427
$ws.Range("N16").Value = -1890.1428
$ws.Range("H31").Value = 33801.49
$ws.Range("J31").Value = 67440.78
$ws.Range("L31").Value = 67440.78
$ws.Range("N31").Value = -68030.78
$ws.Range("H34").Value = 33801.49
$ws.Range("J34").Value = 67440.78
$ws.Range("L34").Value = 67440.78
$ws.Range("N34").Value = -67844.78
$ws.Range("H58").Value = 2053
$ws.Range("I58").Value = 1925.4736
$ws.Range("J58").Value = 2322.2222
$ws.Range("K58").Value = 1925.4736
$ws.Range("L58").Value = 2322.2222
$ws.Range("M58").Value = -1722.4736
$ws.Range("N58").Value = -2728.2222
$ws.Range("H105").Value = 1479.8
$ws.Range("I105").Value = 1599.75
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1599.75
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 147.25
$ws.Range("N105").Value = -4494
$ws.Range("H113").Value = 935.8421
$ws.Range("I113").Value = 714
$ws.Range("J113").Value = 1316.1428
$ws.Range("K113").Value = 714
$ws.Range("L113").Value = 1316.1428
$ws.Range("M113").Value = 1456
$ws.Range("N113").Value = -5656.1428
$ws.Range("H132").Value = 4574.8
$ws.Range("I132").Value = 5728
$ws.Range("J132").Value = 3256.8572
$ws.Range("K132").Value = 17184
$ws.Range("L132").Value = 9770.571599999999
$ws.Range("M132").Value = -14654
$ws.Range("N132").Value = -14830.5716
$ws.Range("H134").Value = 1231.1
$ws.Range("I134").Value = 720.8421
$ws.Range("J134").Value = 2112.4546
$ws.Range("K134").Value = 2162.5263
$ws.Range("L134").Value = 6337.3638
$ws.Range("M134").Value = 372.4737
$ws.Range("N134").Value = -11407.3638
$ws.Range("H136").Value = 2053
$ws.Range("I136").Value = 1925.4736
$ws.Range("J136").Value = 2322.2222
$ws.Range("K136").Value = 5776.4208
$ws.Range("L136").Value = 6966.6666
$ws.Range("M136").Value = -3226.4208
$ws.Range("N136").Value = -12066.6666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15.461538
$ws.Range("I2").Value = 16.5
$ws.Range("J2").Value = 14.571428
$ws.Range("K2").Value = 99
$ws.Range("L2").Value = 87.428568
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = -313.428568
$ws.Range("H5").Value = 1559.7778
$ws.Range("I5").Value = 1340.2222
$ws.Range("J5").Value = 1998.8889
$ws.Range("K5").Value = 4020.6666
$ws.Range("L5").Value = 5996.6667
$ws.Range("M5").Value = -3908.6666
$ws.Range("N5").Value = -6220.6667
$ws.Range("H131").Value = 676135.0600000001
$ws.Range("I131").Value = 675
$ws.Range("J131").Value = 756787
$ws.Range("K131").Value = 2025
$ws.Range("L131").Value = 2270361
$ws.Range("M131").Value = 3015
$ws.Range("N131").Value = -2280441
$ws.Range("H135").Value = 1559.7778
$ws.Range("I135").Value = 1340.2222
$ws.Range("J135").Value = 1998.8889
$ws.Range("K135").Value = 12061.9998
$ws.Range("L135").Value = 17990.0001
$ws.Range("M135").Value = -9526.9998
$ws.Range("N135").Value = -23060.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 72713.73
$ws.Range("I70").Value = 121207.586
$ws.Range("J70").Value = 4014.0833
$ws.Range("K70").Value = 121207.586
$ws.Range("L70").Value = 4014.0833
$ws.Range("M70").Value = -120937.586
$ws.Range("N70").Value = -4554.0833
$ws.Range("H73").Value = 72713.73
$ws.Range("I73").Value = 121207.586
$ws.Range("J73").Value = 4014.0833
$ws.Range("K73").Value = 121207.586
$ws.Range("L73").Value = 4014.0833
$ws.Range("M73").Value = -120271.586
$ws.Range("N73").Value = -5886.0833
$ws.Range("H122").Value = 2237.15
$ws.Range("I122").Value = 1941.6428
$ws.Range("K122").Value = 5824.928400000001
$ws.Range("M122").Value = -3374.928400000001
$ws.Range("H126").Value = 5698.4
$ws.Range("I126").Value = 5623
$ws.Range("K126").Value = 16869
$ws.Range("M126").Value = -14399
$ws.Range("H132").Value = 2248.4138
$ws.Range("I132").Value = 1677.0465
$ws.Range("J132").Value = 3886.3333
$ws.Range("K132").Value = 5031.139499999999
$ws.Range("L132").Value = 11658.9999
$ws.Range("M132").Value = -2501.139499999999
$ws.Range("N132").Value = -16718.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1759.75
$ws.Range("I68").Value = 1787.4286
$ws.Range("J68").Value = 1748.3529
$ws.Range("K68").Value = 1787.4286
$ws.Range("L68").Value = 1748.3529
$ws.Range("M68").Value = -1038.4286
$ws.Range("N68").Value = -3246.3529
$ws.Range("H71").Value = 1759.75
$ws.Range("I71").Value = 1787.4286
$ws.Range("J71").Value = 1748.3529
$ws.Range("K71").Value = 8937.143
$ws.Range("L71").Value = 8741.764500000001
$ws.Range("M71").Value = -5193.143
$ws.Range("N71").Value = -16229.7645
$ws.Range("H100").Value = 2760
$ws.Range("I100").Value = 2116
$ws.Range("J100").Value = 5980
$ws.Range("K100").Value = 2116
$ws.Range("L100").Value = 5980
$ws.Range("M100").Value = -1575
$ws.Range("N100").Value = -7062
$ws.Range("H132").Value = 3822.4546
$ws.Range("I132").Value = 4716.9473
$ws.Range("J132").Value = 2608.5
$ws.Range("K132").Value = 14150.8419
$ws.Range("L132").Value = 7825.5
$ws.Range("M132").Value = -11620.8419
$ws.Range("N132").Value = -12885.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 35845.46
$ws.Range("J115").Value = 35845.46
$ws.Range("L115").Value = 35845.46
$ws.Range("N115").Value = -38979.46
$ws.Range("H132").Value = 22146.686
$ws.Range("I132").Value = 2856.2354
$ws.Range("J132").Value = 60727.59
$ws.Range("K132").Value = 8568.706200000001
$ws.Range("L132").Value = 182182.77
$ws.Range("M132").Value = -6038.706200000001
$ws.Range("N132").Value = -187242.77
